# Update election results for R. A. AÇORES / LAJES DO PICO (row 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 136
$ws.Range("J2").Value = 534
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 152
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 77
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 55
$ws.Range("T2").Value = 93
$ws.Range("U2").Value = 11
$ws.Range("V2").Value = 896
$ws.Range("X2").Value = 859
$ws.Range("Z2").Value = 9
$ws.Range("AA2").Value = 5
